$d = $word.ActiveDocument

$pairs = @(
    @("338×6=2028", "436×4=1744"),
    @("973×3=2919", "385×6=2310"),
    @("494×5=2470", "765×2=1530"),
    @("920×6=5520", "283×2=566"),
    @("283×3=849",  "319×4=1276"),
    @("571×5=2855", "659×5=3295"),
    @("339×3=1017", "492×6=2952"),
    @("673×9=6057", "238×7=1666"),
    @("135×3=405",  "590×6=3540"),
    @("869×5=4345", "131×5=655"),
    @("893×6=5358", "367×2=734"),
    @("383×9=3447", "819×5=4095"),
    @("447×8=3576", "719×4=2876"),
    @("805×5=4025", "182×6=1092"),
    @("197×6=1182", "564×9=5076"),
    @("490×6=2940", "462×9=4158"),
    @("890×5=4450", "381×4=1524"),
    @("293×3=879",  "189×2=378"),
    @("836×5=4180", "563×5=2815"),
    @("775×3=2325", "345×2=690"),
    @("714×9=6426", "892×8=7136"),
    @("657×9=5913", "136×2=272"),
    @("904×5=4520", "516×5=2580"),
    @("870×4=3480", "943×2=1886"),
    @("140×6=840",  "801×8=6408")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

Write-Host "Done replacing $($pairs.Count) values"
